$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.640.53"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3
$ws.Range("D3").Value = "2.111.83"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "347.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5266"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.46%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4509"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.41%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09009"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.44%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.170"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.43%  "

# Row 13
$ws.Range("D13").Value = "2.108.98"
$ws.Range("E13").Value = "  +0.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.804"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.35%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.042"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.81%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001178"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.40%  "

# Row 18
$ws.Range("E18").Value = "  +1.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06713"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.44%  "

# Row 20
$ws.Range("E20").Value = "  +0.37%  "

# Row 21
$ws.Range("E21").Value = "  +1.02%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.300"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.27%  "

# Row 23
$ws.Range("D23").Value = "30.692.17"
$ws.Range("E23").Value = "  +0.63%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.56%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.389"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.06%  "

# Row 26
$ws.Range("D26").Value = "2.361.33"
$ws.Range("E26").Value = "  +0.62%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.41%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.525"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.55%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.88%  "

# Row 31
$ws.Range("E31").Value = "  -0.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1073"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.628"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.42%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.339"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.11%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.017"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.48%  "

# Row 36
$ws.Range("E36").Value = "  -2.24%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.894"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.41%  "

# Row 38
$ws.Range("E38").Value = "  +2.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06836"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.35%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2312"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.32%  "

# Row 41
$ws.Range("E41").Value = "  -1.36%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6870"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.268"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.78%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.13%  "

# Row 45
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.318"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6420"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.61%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.81%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000362"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.50%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.251"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.47%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07289"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.61%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.57%  "
